$wb = $excel.ActiveWorkbook

# Scheduled-runner profit refresh: update currentAveragePrice/NQ/HQ, LevePrice
# NQ/HQ and derived LeveProfit NQ/HQ cells across all job sheets with
# freshly pulled Universalis price data.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 688.9726000000001
$ws.Range("I15").Value = 688.9726000000001
$ws.Range("K15").Value = 2066.9178
$ws.Range("M15").Value = -1897.9178
$ws.Range("H40").Value = 2867.4211
$ws.Range("I40").Value = 2165.4285
$ws.Range("J40").Value = 3276.9167
$ws.Range("K40").Value = 2165.4285
$ws.Range("L40").Value = 3276.9167
$ws.Range("M40").Value = -1990.4285
$ws.Range("N40").Value = -3626.9167
$ws.Range("H51").Value = 3949.5
$ws.Range("I51").Value = 1898
$ws.Range("J51").Value = 4359.8
$ws.Range("K51").Value = 1898
$ws.Range("L51").Value = 4359.8
$ws.Range("M51").Value = -1414
$ws.Range("N51").Value = -5327.8
$ws.Range("H70").Value = 2281
$ws.Range("J70").Value = 2528.8333
$ws.Range("L70").Value = 7586.499899999999
$ws.Range("N70").Value = -8126.499899999999
$ws.Range("H73").Value = 2281
$ws.Range("J73").Value = 2528.8333
$ws.Range("L73").Value = 7586.499899999999
$ws.Range("N73").Value = -9458.499899999999
$ws.Range("H76").Value = 11768.941
$ws.Range("I76").Value = 15730.889
$ws.Range("K76").Value = 15730.889
$ws.Range("M76").Value = -15415.889
$ws.Range("H79").Value = 11768.941
$ws.Range("I79").Value = 15730.889
$ws.Range("K79").Value = 15730.889
$ws.Range("M79").Value = -14638.889
$ws.Range("H96").Value = 3216.6072
$ws.Range("I96").Value = 2366.6191
$ws.Range("J96").Value = 5766.5713
$ws.Range("K96").Value = 7099.8573
$ws.Range("L96").Value = 17299.7139
$ws.Range("M96").Value = -5726.8573
$ws.Range("N96").Value = -20045.7139
$ws.Range("H116").Value = 3633.5557
$ws.Range("I116").Value = 3783.3
$ws.Range("K116").Value = 3783.3
$ws.Range("M116").Value = -341.3000000000002
$ws.Range("H121").Value = 4616.6665
$ws.Range("J121").Value = 4616.6665
$ws.Range("L121").Value = 13849.9995
$ws.Range("N121").Value = -17343.9995
$ws.Range("H137").Value = 2370.6
$ws.Range("I137").Value = 2363.3
$ws.Range("K137").Value = 7089.900000000001
$ws.Range("M137").Value = -4539.900000000001
$ws.Range("H141").Value = 2958.5789
$ws.Range("I141").Value = 2956.3333
$ws.Range("K141").Value = 8868.999899999999
$ws.Range("M141").Value = -3688.999899999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 27781540
$ws.Range("I61").Value = 35716588
$ws.Range("J61").Value = 8874.875
$ws.Range("K61").Value = 35716588
$ws.Range("L61").Value = 8874.875
$ws.Range("M61").Value = -35716376
$ws.Range("N61").Value = -9298.875
$ws.Range("H122").Value = 2992.1292
$ws.Range("I122").Value = 1285.2142
$ws.Range("K122").Value = 3855.6426
$ws.Range("M122").Value = -1405.6426
$ws.Range("H136").Value = 27781540
$ws.Range("I136").Value = 35716588
$ws.Range("J136").Value = 8874.875
$ws.Range("K136").Value = 107149764
$ws.Range("L136").Value = 26624.625
$ws.Range("M136").Value = -107147214
$ws.Range("N136").Value = -31724.625
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 52977.832
$ws.Range("I86").Value = 27966.75
$ws.Range("K86").Value = 27966.75
$ws.Range("M86").Value = -26843.75
$ws.Range("H89").Value = 52977.832
$ws.Range("I89").Value = 27966.75
$ws.Range("K89").Value = 139833.75
$ws.Range("M89").Value = -134217.75
$ws.Range("H105").Value = 11183.8
$ws.Range("I105").Value = 17640.5
$ws.Range("J105").Value = 1498.75
$ws.Range("K105").Value = 17640.5
$ws.Range("L105").Value = 1498.75
$ws.Range("M105").Value = -15893.5
$ws.Range("N105").Value = -4992.75
$ws.Range("H107").Value = 3663.6843
$ws.Range("I107").Value = 3289.5386
$ws.Range("K107").Value = 3289.5386
$ws.Range("M107").Value = -1369.5386
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 18524490
$ws.Range("I31").Value = 5354.9
$ws.Range("J31").Value = 71436300
$ws.Range("K31").Value = 5354.9
$ws.Range("L31").Value = 71436300
$ws.Range("M31").Value = -5059.9
$ws.Range("N31").Value = -71436890
$ws.Range("H34").Value = 18524490
$ws.Range("I34").Value = 5354.9
$ws.Range("J34").Value = 71436300
$ws.Range("K34").Value = 5354.9
$ws.Range("L34").Value = 71436300
$ws.Range("M34").Value = -5152.9
$ws.Range("N34").Value = -71436704
$ws.Range("H58").Value = 2710.2
$ws.Range("I58").Value = 2344.9048
$ws.Range("J58").Value = 3562.5557
$ws.Range("K58").Value = 2344.9048
$ws.Range("L58").Value = 3562.5557
$ws.Range("M58").Value = -2141.9048
$ws.Range("N58").Value = -3968.5557
$ws.Range("H105").Value = 8196.0625
$ws.Range("I105").Value = 1642.6666
$ws.Range("J105").Value = 16621.857
$ws.Range("K105").Value = 1642.6666
$ws.Range("L105").Value = 16621.857
$ws.Range("M105").Value = 104.3334
$ws.Range("N105").Value = -20115.857
$ws.Range("H132").Value = 53736
$ws.Range("I132").Value = 70593.39999999999
$ws.Range("K132").Value = 211780.2
$ws.Range("M132").Value = -209250.2
$ws.Range("H133").Value = 326163
$ws.Range("J133").Value = 326163
$ws.Range("L133").Value = 326163
$ws.Range("N133").Value = -331223
$ws.Range("H134").Value = 1964.8182
$ws.Range("I134").Value = 1761.3
$ws.Range("K134").Value = 5283.9
$ws.Range("M134").Value = -2748.9
$ws.Range("H136").Value = 2710.2
$ws.Range("I136").Value = 2344.9048
$ws.Range("J136").Value = 3562.5557
$ws.Range("K136").Value = 7034.714399999999
$ws.Range("L136").Value = 10687.6671
$ws.Range("M136").Value = -4484.714399999999
$ws.Range("N136").Value = -15787.6671
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3563.5
$ws.Range("I3").Value = 795.5714
$ws.Range("J3").Value = 10022
$ws.Range("K3").Value = 2386.7142
$ws.Range("L3").Value = 30066
$ws.Range("M3").Value = -2274.7142
$ws.Range("N3").Value = -30290
$ws.Range("H107").Value = 1215.4166
$ws.Range("J107").Value = 1880.5714
$ws.Range("L107").Value = 5641.7142
$ws.Range("N107").Value = -9481.7142
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 380086
$ws.Range("I33").Value = 6450
$ws.Range("J33").Value = 504631.34
$ws.Range("K33").Value = 6450
$ws.Range("L33").Value = 504631.34
$ws.Range("M33").Value = -6198
$ws.Range("N33").Value = -505135.34
$ws.Range("H41").Value = 2117.25
$ws.Range("I41").Value = 2117.25
$ws.Range("K41").Value = 2117.25
$ws.Range("M41").Value = -1762.25
$ws.Range("H80").Value = 3254.6667
$ws.Range("J80").Value = 3533.2856
$ws.Range("L80").Value = 3533.2856
$ws.Range("N80").Value = -5529.2856
$ws.Range("H83").Value = 3254.6667
$ws.Range("J83").Value = 3533.2856
$ws.Range("L83").Value = 17666.428
$ws.Range("N83").Value = -27650.428
$ws.Range("H102").Value = 2816.2727
$ws.Range("I102").Value = 1938.9231
$ws.Range("J102").Value = 4083.5557
$ws.Range("K102").Value = 1938.9231
$ws.Range("L102").Value = 4083.5557
$ws.Range("M102").Value = -316.9231
$ws.Range("N102").Value = -7327.5557
$ws.Range("H122").Value = 2275.238
$ws.Range("I122").Value = 2043.3334
$ws.Range("J122").Value = 3666.6667
$ws.Range("K122").Value = 6130.0002
$ws.Range("L122").Value = 11000.0001
$ws.Range("M122").Value = -3680.0002
$ws.Range("N122").Value = -15900.0001
$ws.Range("H123").Value = 55714
$ws.Range("J123").Value = 70000
$ws.Range("L123").Value = 70000
$ws.Range("N123").Value = -74900
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4132.383
$ws.Range("I7").Value = 3613.762
$ws.Range("J7").Value = 4551.269
$ws.Range("K7").Value = 3613.762
$ws.Range("L7").Value = 4551.269
$ws.Range("M7").Value = -3501.762
$ws.Range("N7").Value = -4775.269
$ws.Range("H126").Value = 4132.383
$ws.Range("I126").Value = 3613.762
$ws.Range("J126").Value = 4551.269
$ws.Range("K126").Value = 10841.286
$ws.Range("L126").Value = 13653.807
$ws.Range("M126").Value = -8371.286
$ws.Range("N126").Value = -18593.807
$ws.Range("H132").Value = 133335360
$ws.Range("I132").Value = 1900.4286
$ws.Range("J132").Value = 250002140
$ws.Range("K132").Value = 5701.2858
$ws.Range("L132").Value = 750006420
$ws.Range("M132").Value = -3171.2858
$ws.Range("N132").Value = -750011480
$ws.Range("H136").Value = 2444.7415
$ws.Range("I136").Value = 2090.1177
$ws.Range("K136").Value = 6270.353099999999
$ws.Range("M136").Value = -3720.353099999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 38504576
$ws.Range("I122").Value = 43526576
$ws.Range("K122").Value = 130579728
$ws.Range("M122").Value = -130577278
$ws.Range("H126").Value = 4876
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 5684.718
$ws.Range("I132").Value = 5650.3687
$ws.Range("K132").Value = 16951.1061
$ws.Range("M132").Value = -14421.1061
